$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.173.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.82%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.653.90'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.55%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  +0.33%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''217.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +0.04%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.5312'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +1.30%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.32%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.2622'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +0.22%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.06331'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +0.93%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  +0.48%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.07802'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +0.77%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''4.519'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +1.48%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''1.637.02'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -0.32%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''1.881.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +0.58%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.5489'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +1.04%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.0₅8183'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +1.42%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''65.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +0.92%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''26.152.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +0.64%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''1.006'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +0.26%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''4.596'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +1.15%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''191.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -0.36%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''10.09'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +0.58%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''6.008'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +0.57%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''  +0.35%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''145.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +3.96%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  -0.88%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''7.197'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -0.68%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -1.04%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''1.472'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +3.64%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''0.05737'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -3.38%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''1.272'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.03%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  +1.85%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''3.261'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +0.91%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''1.593'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +4.71%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''2.807'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +2.23%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''2.422'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +0.39%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''0.9485'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +1.00%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.5741'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +0.59%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.01602'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -0.05%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = '''TrustWalletToken'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = '''0.8491'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +0.34%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = '''FraxShare'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = '''5.783'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -1.19%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("B43").Value = '''Quant'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = '''https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = '''103.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +3.38%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = '''Maker'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = '''1.038.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +3.54%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''1.793.62'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +0.43%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''56.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +0.39%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = '''BabyDogeCoin'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''0.0₈105'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +0.08%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = '''Frax'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = '''1.005'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -0.05%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.4357'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +1.69%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.05156'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +0.15%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''7.847'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.16%  '
$ws.Range("E51").Style = "Normal"
